$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 56 -- this shifts the existing rows 56-86 down to 57-87,
# matching the diff (dimension grows from A1:R86 to A1:R87).
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record.
$ws.Range("A56").Value = 11
$ws.Range("B56").Value = "Vega Monumental Concepción"
$ws.Range("C56").Value = "Bíobío"
$ws.Range("D56").Value = 45119
$ws.Range("E56").Value = 8
$ws.Range("F56").Value = 100112031
$ws.Range("G56").Value = "Poroto verde"
$ws.Range("H56").Value = "Magnum"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 100
$ws.Range("K56").Value = 25000
$ws.Range("L56").Value = 26000
$ws.Range("M56").Value = 25500
$ws.Range("N56").Value = "$/malla 25 kilos"
$ws.Range("O56").Value = "Perú"
$ws.Range("P56").Value = 1020
$ws.Range("Q56").Value = 25
$ws.Range("R56").Value = "Hortaliza"
